$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised AgTests / AgPosit figures for rows 595-627 (retroactive corrections)
$ws.Range("F595").Value = 27268
$ws.Range("F597").Value = 29581
$ws.Range("F600").Value = 40557
$ws.Range("F602").Value = 30087
$ws.Range("F603").Value = 32025
$ws.Range("F604").Value = 29963
$ws.Range("F605").Value = 14803
$ws.Range("F606").Value = 14328
$ws.Range("G606").Value = 1264
$ws.Range("F608").Value = 46038
$ws.Range("G608").Value = 2896
$ws.Range("F609").Value = 36406
$ws.Range("G609").Value = 2162
$ws.Range("F610").Value = 33922
$ws.Range("F611").Value = 34102
$ws.Range("G611").Value = 2126
$ws.Range("F612").Value = 16282
$ws.Range("F613").Value = 21605
$ws.Range("F614").Value = 47687
$ws.Range("G614").Value = 3328
$ws.Range("F615").Value = 36682
$ws.Range("G615").Value = 2350
$ws.Range("F616").Value = 38151
$ws.Range("G616").Value = 2563
$ws.Range("F617").Value = 38838
$ws.Range("G617").Value = 2592
$ws.Range("F618").Value = 37616
$ws.Range("F619").Value = 17704
$ws.Range("G619").Value = 1883
$ws.Range("F620").Value = 25516
$ws.Range("G620").Value = 2373
$ws.Range("F621").Value = 55537
$ws.Range("G621").Value = 4083
$ws.Range("F622").Value = 40870
$ws.Range("G622").Value = 2980
$ws.Range("F623").Value = 14747
$ws.Range("G623").Value = 1539
$ws.Range("F624").Value = 50575
$ws.Range("G624").Value = 3923
$ws.Range("F625").Value = 43130
$ws.Range("G625").Value = 3515
$ws.Range("F626").Value = 19708
$ws.Range("G626").Value = 2066
$ws.Range("F627").Value = 32813
$ws.Range("G627").Value = 2663

# New daily rows appended for 2021-11-22, 2021-11-23, 2021-11-24
$ws.Range("A628").Value = 44522
$ws.Range("A628").NumberFormat = "yyyy-mm-dd"
$ws.Range("B628").Value = 621423
$ws.Range("C628").Value = 23107
$ws.Range("D628").Value = 6739
$ws.Range("E628").Value = 13985
$ws.Range("F628").Value = 62422
$ws.Range("G628").Value = 4078

$ws.Range("A629").Value = 44523
$ws.Range("A629").NumberFormat = "yyyy-mm-dd"
$ws.Range("B629").Value = 631738
$ws.Range("C629").Value = 30175
$ws.Range("D629").Value = 10315
$ws.Range("E629").Value = 14056
$ws.Range("F629").Value = 44018
$ws.Range("G629").Value = 2772

$ws.Range("A630").Value = 44524
$ws.Range("A630").NumberFormat = "yyyy-mm-dd"
$ws.Range("B630").Value = 641903
$ws.Range("C630").Value = 30620
$ws.Range("D630").Value = 10165
$ws.Range("E630").Value = 14107
$ws.Range("F630").Value = 33819
$ws.Range("G630").Value = 2063
